# Insert a new daily price record as row 6 on Sheet1, pushing the existing
# rows 6-84 down to 7-85 (dimension grows from A1:R84 to A1:R85).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6..84 down by one, leaving a blank row 6 (with D6 keeping the
# date-column style carried over from the former row 6).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new "Papa" price entry.
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44817
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100114001
$ws.Range("G6").Value = "Papa"
$ws.Range("H6").Value = "Asterix"
$ws.Range("I6").Value = "1a (cosecha lavada)"
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 13000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 13500
$ws.Range("N6").Value = "$/malla 25 kilos"
$ws.Range("O6").Value = "Región de Los Lagos"
$ws.Range("P6").Value = 540
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
